$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The text-valued columns (A, B, D, E) hold numeric-looking strings that must
# stay text (shared-string) cells rather than be auto-coerced to numbers, so
# format them as Text first. New values are written column-by-column
# (A2,A3,B2,B3,D2,D3,E2,E3) so the shared-string table is rebuilt in the
# same relative order the source workbook used.
$ws.Range("A2:B3").NumberFormat = "@"
$ws.Range("D2:E3").NumberFormat = "@"

$ws.Range("A2").Value = "8917.1622553"
$ws.Range("A3").Value = "8981.6220102"

$ws.Range("B2").Value = "8918.401366"
$ws.Range("B3").Value = "8981.656001"

$ws.Range("D2").Value = "8930.6574626"
$ws.Range("D3").Value = "8995.0509062"

$ws.Range("E2").Value = "8967.0377007"
$ws.Range("E3").Value = "9031.1469414"

# Column C switches from text to plain numbers.
$ws.Range("C2").Value = 144
$ws.Range("C3").Value = 307.2

# Column F numeric values change.
$ws.Range("F2").Value = 1.23911070000031
$ws.Range("F3").Value = 0.03399079999871901
